# Weekly fruit/vegetable price update: a new weekly record is inserted for
# this market/variety at row 57 (pushing the existing rows 57-132 down to
# 58-133, which keeps all of their data unchanged), and the new row is
# populated with this week's reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 57, shifting existing rows 57:132 down to 58:133.
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new weekly record.
$ws.Cells.Item(57, 1).Value = 4
$ws.Cells.Item(57, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(57, 3).Value = "Los Lagos"
$ws.Cells.Item(57, 4).Value = "2021-09-28"
$ws.Cells.Item(57, 5).Value = 10
$ws.Cells.Item(57, 6).Value = 100112017
$ws.Cells.Item(57, 7).Value = "Apio"
$ws.Cells.Item(57, 8).Value = "Americana (o)"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 40
$ws.Cells.Item(57, 11).Value = 12000
$ws.Cells.Item(57, 12).Value = 12000
$ws.Cells.Item(57, 13).Value = 12000
$ws.Cells.Item(57, 14).Value = "`$/docena de matas"
$ws.Cells.Item(57, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(57, 16).Value = 2000
$ws.Cells.Item(57, 17).Value = 6
$ws.Cells.Item(57, 18).Value = "Hortaliza"
